# Insert a new weekly price record (row 57) for Locoto - Primera,
# pushing the existing rows 57:141 down to 58:142.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(57).Insert()

$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(57, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value = 44848
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(57, 6).Value = 100112042
$ws.Cells.Item(57, 7).Value = "Locoto"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 160
$ws.Cells.Item(57, 11).Value = 18000
$ws.Cells.Item(57, 12).Value = 19000
$ws.Cells.Item(57, 13).Value = 18500
$ws.Cells.Item(57, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(57, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 16).Value = 925
$ws.Cells.Item(57, 17).Value = 20
$ws.Cells.Item(57, 18).Value = "Hortaliza"
